$d = $word.ActiveDocument

# --- 1. Update the placeholder ID text in the first paragraph, and drop the
#        trailing single-space run that followed it. ---
$d.Content.Find.Execute("**ID__AFFARS_5333_topic_12__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5333_290__ID**", 2)

$p1 = $d.Paragraphs(1)

# Remove the trailing " " run left over after the placeholder text (the
# paragraph now reads "**ID__AFFARS_5333_290__ID** " - strip that trailing
# space so only the single run with the new id text remains).
$pRange = $p1.Range
$pRange.End = $pRange.End - 1
$tail = $d.Range($pRange.End - 1, $pRange.End)
if ($tail.Text -eq " ") {
    $tail.Text = ""
}

# --- 2. Paragraph formatting: indent + paragraph border ---
$p1.Range.ParagraphFormat.LeftIndent = 11.25

$borders = $p1.Range.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5
